$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old second row (the single date value row) -- the new layout
# only has a header row.
$ws.Range("A2").EntireRow.Delete()

# Populate the header row with the list of bluesky follower handles.
$ws.Cells.Item(1,1).Value = "date"
$ws.Cells.Item(1,2).Value = "nea-ukcharity.bsky.social"
$ws.Cells.Item(1,3).Value = "caneurope.bsky.social"
$ws.Cells.Item(1,4).Value = "wateraid.bsky.social"
$ws.Cells.Item(1,5).Value = "migrantsrights.bsky.social"
$ws.Cells.Item(1,6).Value = "friends-earth.bsky.social"
$ws.Cells.Item(1,7).Value = "samcardwell44.bsky.social"
$ws.Cells.Item(1,8).Value = "greenpeace.eu"
$ws.Cells.Item(1,9).Value = "wwfeu.bsky.social"
$ws.Cells.Item(1,10).Value = "powertochange.org.uk"
$ws.Cells.Item(1,11).Value = "thegreenregister.bsky.social"
$ws.Cells.Item(1,12).Value = "endfuelpoverty.bsky.social"
$ws.Cells.Item(1,13).Value = "commenergyengland.bsky.social"
$ws.Cells.Item(1,14).Value = "extinctionrebellion.uk"
$ws.Cells.Item(1,15).Value = "wwtworldwide.bsky.social"
$ws.Cells.Item(1,16).Value = "bristolgreenparty.bsky.social"
$ws.Cells.Item(1,17).Value = "warmthiswinter.bsky.social"
$ws.Cells.Item(1,18).Value = "jrct.bsky.social"
$ws.Cells.Item(1,19).Value = "ssencommunity.bsky.social"
$ws.Cells.Item(1,20).Value = "localtrust.bsky.social"
$ws.Cells.Item(1,21).Value = "wiltscouncil.bsky.social"
$ws.Cells.Item(1,22).Value = "nationalgrid.bsky.social"
$ws.Cells.Item(1,23).Value = "ofgem.bsky.social"
$ws.Cells.Item(1,24).Value = "barnsleycouncil.bsky.social"
$ws.Cells.Item(1,25).Value = "northsomersetc.bsky.social"
$ws.Cells.Item(1,26).Value = "citizensadvice.bsky.social"

$ws.Range("A1:Z1").Select()
